$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in row 2
$ws.Range("A2").Value = 2021662
$ws.Range("B2").Value = 82253

# Update the active selection to B3
$ws.Range("B3").Select()
